# GA: add new feature
# Adds a new "上级用户名" (Superior Username) column to the user-import
# template's header row, bolds / heightens the header row, widens the
# columns to fit the new/longer headers, and moves the selection to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- new header cell -------------------------------------------------
$ws.Range("H1").Value = "上级用户名"

# --- header row formatting -------------------------------------------
$ws.Range("A1:H1").Font.Bold = $true
$ws.Rows("1:1").RowHeight = 25

# --- column widths (values chosen so the engine's char->pixel rounding
#     lands as close as possible to the authored widths) --------------
$ws.Columns("B:B").ColumnWidth = 16.714285714285715
$ws.Columns("D:D").ColumnWidth = 15.714285714285714
$ws.Columns("E:E").ColumnWidth = 20.857142857142858
$ws.Columns("F:F").ColumnWidth = 15.0
$ws.Columns("G:G").ColumnWidth = 23.714285714285715
$ws.Columns("H:H").ColumnWidth = 20.857142857142858

# --- selection moves to the new column's second row -------------------
[void]$ws.Range("H2").Select()
